$d = $word.ActiveDocument

# 1. Trim "from left to right" out of the intro sentence.
$d.Content.Find.Execute(
    "The buttons at the bottom from left to right, with their keyboard equivalents are ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The buttons at the bottom with their keyboard equivalents are ", 2)

# 2. Add a new "Type / T / Change between D and G Merlin" row below the
#    "Start / G / Resume" row, moving the _GoBack bookmark down onto it.

# The hidden _GoBack bookmark currently sits at the end of the "Resume" row;
# remove it from there so it can be re-created at the end of the new row.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Locate the "Start ... Resume" paragraph and append a fresh paragraph after it.
$resumeParagraph = $d.Paragraphs.Item(8)
$resumeParagraph.Range.InsertParagraphAfter()
$newParagraph = $d.Paragraphs.Item(9)

# Build the new row with real <w:tab/> runs (matching the other rows) and
# re-attach the _GoBack bookmark at its end via a raw OOXML fragment, since
# that is the only reliable way to get literal tab runs + the bookmark pair
# in one shot.
$newRowXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">Type </w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t>T</w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t>Change between D and G Merlin</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$newParagraph.Range.InsertXML($newRowXml)
